# Slide 23 ("Final Analysis") -> Content Placeholder 2 -> 3rd paragraph
# (the ANOVA summary bullet). Replace the second half of the sentence
# ("...This demonstrates that the Reading scores ... significant
# difference.") with the new conclusion, delivered as two runs so the
# trailing clause "the regions." lands in its own run, matching the
# authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(3)

$run1Text = "The Analysis of Variance Test (ANOVA) was performed on the Math and Reading Scores for both 4th grade and 8th grade across the regions. It is observed that the p-values that were computed from the ANOVA test for Math Grade 4 and Grade 8 comes out to be 6.76 e -05 and 2.07 e -05 where as the Reading Grade 4 and Grade 8 p-value scores are 0.0001 and 0.57 across regions.  This demonstrates that the Math scores are more statistically significant than the reading scores among "
$run2Text = "the regions."

# Clear the paragraph first so the upcoming assignment isn't diffed
# against the old wording (which would otherwise let the host reuse
# fragments of the previous sentence as extra runs).
$para.Text = "#"

$tr.Paragraphs(3).Text = $run1Text
$tr.Paragraphs(3).InsertAfter($run2Text) | Out-Null
